$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text representation (e.g. "1.00", "0.990")
# instead of being auto-converted to numbers by Excel when set via COM.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '51.352.97'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '2.980.31'
$ws.Range('E3').Value = '  +2.12%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '383.80'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D6').Value = '102.83'
$ws.Range('E6').Value = '  +0.66%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.592'
$ws.Range('E9').Value = '  +1.39%  '
$ws.Range('D10').Value = '36.71'
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').Value = '0.0842'
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('D13').Value = '3.450.47'
$ws.Range('E13').Value = '  +1.96%  '
$ws.Range('D14').Value = '18.17'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('D15').Value = '7.50'
$ws.Range('E15').Value = '  +2.33%  '
$ws.Range('D16').Value = '2.984.04'
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('D17').Value = '0.990'
$ws.Range('E17').Value = '  +7.04%  '
$ws.Range('D18').Value = '51.328.02'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = '3.25'
$ws.Range('E19').Value = '  -4.84%  '
$ws.Range('D20').Value = '7.39'
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('D21').Value = '12.81'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = '0.0₃0957'
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('D23').Value = '68.80'
$ws.Range('E23').Value = '  +0.94%  '
$ws.Range('D24').Value = '261.71'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').Value = '2.90'
$ws.Range('E25').Value = '  +4.96%  '
$ws.Range('D26').Value = '8.23'
$ws.Range('E26').Value = '  +15.67%  '
$ws.Range('D27').Value = '7.54'
$ws.Range('E27').Value = '  +11.11%  '
$ws.Range('D28').Value = '4.15'
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '0.114'
$ws.Range('E29').Value = '  +12.84%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').Value = '0.168'
$ws.Range('E30').Value = '  -0.84%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').Value = '25.82'
$ws.Range('E32').Value = '  +0.95%  '
$ws.Range('D33').Value = '9.84'
$ws.Range('E33').Value = '  +0.94%  '
$ws.Range('D34').Value = '34.46'
$ws.Range('E34').Value = '  +1.64%  '
$ws.Range('D35').Value = '51.00'
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('D36').Value = '2.06'
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('D37').Value = '0.0449'
$ws.Range('E37').Value = '  +7.19%  '
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('D39').Value = '2.99'
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').Value = '17.04'
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range('D41').Value = '2.58'
$ws.Range('E41').Value = '  +1.48%  '
$ws.Range('D42').Value = '0.115'
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('D43').Value = '1.81'
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('D44').Value = '122.23'
$ws.Range('E44').Value = '  -0.48%  '
$ws.Range('D45').Value = '21.70'
$ws.Range('E45').Value = '  +1.03%  '
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '0.275'
$ws.Range('E47').Value = '  +2.64%  '
$ws.Range('D48').Value = '2.34'
$ws.Range('E48').Value = '  +1.20%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.024.98'
$ws.Range('E49').Value = '  -0.10%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').Value = '3.24'
$ws.Range('E50').Value = '  +3.43%  '
$ws.Range('D51').Value = '0.0329'
$ws.Range('E51').Value = '  +3.43%  '
